$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1, reusing the same formatting as the existing headers (A1:E1)
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Fill boolean outlier flags for rows 2-25 in columns F, G, H (all FALSE by default)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}

# Row 24 has KNN_Outliers_MAD = TRUE
$ws.Cells.Item(24, 6).Value = $true
